$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26 (shifts existing rows 26-79 down to 27-80)
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with this week's new record
$ws.Range("A26").Value = 6
$ws.Range("B26").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44614
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100101
$ws.Range("H26").Value = "Berries"
$ws.Range("I26").Value = 100101008
$ws.Range("J26").Value = "Mora"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 6000
$ws.Range("O26").Value = 6000
$ws.Range("P26").Value = 6000
$ws.Range("Q26").Value = "$/bandeja 2 kilos"
$ws.Range("R26").Value = "Provincia de Curicó"
$ws.Range("S26").Value = 3000
$ws.Range("T26").Value = 2

# Insert a second new row at row 78 (shifts current rows 78-80 down to 79-81)
$ws.Rows.Item(78).Insert()

# Populate the new row 78 with the other new record for this week
$ws.Range("A78").Value = 6
$ws.Range("B78").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C78").Value = "Metropolitana"
$ws.Range("D78").Value = 44615
$ws.Range("E78").Value = 13
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100101
$ws.Range("H78").Value = "Berries"
$ws.Range("I78").Value = 100101008
$ws.Range("J78").Value = "Mora"
$ws.Range("K78").Value = "Sin especificar"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 50
$ws.Range("N78").Value = 6000
$ws.Range("O78").Value = 6000
$ws.Range("P78").Value = 6000
$ws.Range("Q78").Value = "$/bandeja 2 kilos"
$ws.Range("R78").Value = "Provincia de Curicó"
$ws.Range("S78").Value = 3000
$ws.Range("T78").Value = 2
